# "Update countries & provincias Spain"
#
# Refreshes the COVID country table on sheet "Pais" with newer figures and
# re-sorts a few rows whose "Casos totales" (column B) now rank differently:
#   - Moldavia overtakes Costa Rica (rows 65/66)
#   - Eslovenia overtakes Lituania (rows 129/130)
#   - Santa Lucia overtakes Timor Oriental (rows 202/203)
# Also bumps the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Last updated" timestamp -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 16:56"

# --- Row 4: Estados Unidos ------------------------------------------------
$ws.Cells.Item(4, 2).Value = 5847413
$ws.Cells.Item(4, 3).Value = 5985
$ws.Cells.Item(4, 4).Value = 3148159
$ws.Cells.Item(4, 5).Value = 2518976
$ws.Cells.Item(4, 7).Value = 104
$ws.Cells.Item(4, 8).Value = 180278

# --- Row 15: Argentina ----------------------------------------------------
$ws.Cells.Item(15, 4).Value = 251400
$ws.Cells.Item(15, 5).Value = 78455
$ws.Cells.Item(15, 7).Value = 99
$ws.Cells.Item(15, 8).Value = 6947

# --- Row 23: Alemania ------------------------------------------------------
$ws.Cells.Item(23, 2).Value = 234151
$ws.Cells.Item(23, 3).Value = 294
$ws.Cells.Item(23, 5).Value = 15869
$ws.Cells.Item(23, 7).Value = 1
$ws.Cells.Item(23, 8).Value = 9332

# --- Row 50: Portugal -------------------------------------------------------
$ws.Cells.Item(50, 2).Value = 55597
$ws.Cells.Item(50, 3).Value = 145
$ws.Cells.Item(50, 4).Value = 40774
$ws.Cells.Item(50, 5).Value = 13027
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(50, 8).Value = 1796

# --- Row 65: was Costa Rica, now Moldavia (new figures, overtakes row 66) ---
$ws.Cells.Item(65, 1).Value = "Moldavia"
$ws.Cells.Item(65, 2).Value = 33478
$ws.Cells.Item(65, 3).Value = 406
$ws.Cells.Item(65, 4).Value = 22683
$ws.Cells.Item(65, 5).Value = 9855
$ws.Cells.Item(65, 7).Value = 5
$ws.Cells.Item(65, 8).Value = 940

# --- Row 66: was Moldavia, now Costa Rica (keeps its prior figures) ---------
$ws.Cells.Item(66, 1).Value = "Costa Rica"
$ws.Cells.Item(66, 2).Value = 33084
$ws.Cells.Item(66, 4).Value = 10372
$ws.Cells.Item(66, 5).Value = 22364
$ws.Cells.Item(66, 8).Value = 348

# --- Row 67: Kenia -----------------------------------------------------------
$ws.Cells.Item(67, 2).Value = 32364
$ws.Cells.Item(67, 3).Value = 246
$ws.Cells.Item(67, 4).Value = 18670
$ws.Cells.Item(67, 5).Value = 13146
$ws.Cells.Item(67, 7).Value = 6
$ws.Cells.Item(67, 8).Value = 548

# --- Row 73: El Salvador -------------------------------------------------------
$ws.Cells.Item(73, 2).Value = 24622
$ws.Cells.Item(73, 3).Value = 202
$ws.Cells.Item(73, 4).Value = 12246
$ws.Cells.Item(73, 5).Value = 11715

# --- Row 88: Zambia ---------------------------------------------------------
$ws.Cells.Item(88, 2).Value = 11082
$ws.Cells.Item(88, 3).Value = 251
$ws.Cells.Item(88, 5).Value = 860
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = 280

# --- Row 129: was Lituania, now Eslovenia (new figures, overtakes row 130) --
$ws.Cells.Item(129, 1).Value = "Eslovenia"
$ws.Cells.Item(129, 2).Value = 2651
$ws.Cells.Item(129, 3).Value = 34
$ws.Cells.Item(129, 4).Value = 2079
$ws.Cells.Item(129, 5).Value = 441
$ws.Cells.Item(129, 8).Value = 131

# --- Row 130: was Eslovenia, now Lituania (keeps its prior figures) ---------
$ws.Cells.Item(130, 1).Value = "Lituania"
$ws.Cells.Item(130, 2).Value = 2635
$ws.Cells.Item(130, 3).Value = 41
$ws.Cells.Item(130, 4).Value = 1766
$ws.Cells.Item(130, 5).Value = 785
$ws.Cells.Item(130, 8).Value = 84

# --- Row 131: Sudan del Sur ---------------------------------------------------
$ws.Cells.Item(131, 2).Value = 2499
$ws.Cells.Item(131, 3).Value = 2
$ws.Cells.Item(131, 5).Value = 1162

# --- Row 162: Trinidad yTobago -----------------------------------------------
$ws.Cells.Item(162, 2).Value = 963
$ws.Cells.Item(162, 3).Value = 33
$ws.Cells.Item(162, 5).Value = 784
$ws.Cells.Item(162, 7).Value = 1
$ws.Cells.Item(162, 8).Value = 14

# --- Row 171: Birmania ---------------------------------------------------------
$ws.Cells.Item(171, 2).Value = 450
$ws.Cells.Item(171, 3).Value = 9
$ws.Cells.Item(171, 4).Value = 341
$ws.Cells.Item(171, 5).Value = 103

# --- Row 202: was Timor Oriental, now Santa Lucia (figures already tied) -----
$ws.Cells.Item(202, 1).Value = "Santa Lucia"

# --- Row 203: was Santa Lucia, now Timor Oriental (figures already tied) -----
$ws.Cells.Item(203, 1).Value = "Timor Oriental"
